$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $value)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '26.215.17'
Set-TextValue 'D3' '1.659.50'
Set-TextValue 'E3' '  -0.43%  '
Set-TextValue 'D4' '1.004'
Set-TextValue 'E4' '  -0.68%  '
Set-TextValue 'D5' '218.44'
Set-TextValue 'E5' '  -0.40%  '
Set-TextValue 'D6' '0.5238'
Set-TextValue 'E6' '  -2.04%  '
Set-TextValue 'D7' '1.005'
Set-TextValue 'E7' '  -0.59%  '
Set-TextValue 'D8' '0.2641'
Set-TextValue 'E8' '  -0.73%  '
Set-TextValue 'D9' '0.06321'
Set-TextValue 'E9' '  -1.02%  '
Set-TextValue 'D10' '20.70'
Set-TextValue 'E10' '  -0.68%  '
Set-TextValue 'D11' '0.07798'
Set-TextValue 'E11' '  -0.48%  '
Set-TextValue 'B12' 'WrappedEther'
Set-TextValue 'C12' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D12' '1.762.80'
Set-TextValue 'E12' '  +5.72%  '
Set-TextValue 'B13' 'Polkadot'
Set-TextValue 'C13' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D13' '4.508'
Set-TextValue 'E13' '  -1.04%  '
Set-TextValue 'D14' '1.889.51'
Set-TextValue 'E14' '  -0.27%  '
Set-TextValue 'D15' '0.5633'
Set-TextValue 'E15' '  +1.65%  '
Set-TextValue 'D16' '0.0₅8065'
Set-TextValue 'E16' '  -1.55%  '
Set-TextValue 'D17' '65.26'
Set-TextValue 'E17' '  -0.97%  '
Set-TextValue 'D18' '26.227.49'
Set-TextValue 'E18' '  -0.46%  '
Set-TextValue 'E19' '  -0.65%  '
Set-TextValue 'D20' '4.716'
Set-TextValue 'E20' '  +0.98%  '
Set-TextValue 'D21' '194.43'
Set-TextValue 'E21' '  +0.35%  '
Set-TextValue 'D22' '10.23'
Set-TextValue 'E22' '  -0.40%  '
Set-TextValue 'D23' '6.017'
Set-TextValue 'E23' '  -0.37%  '
Set-TextValue 'D24' '1.005'
Set-TextValue 'E24' '  -0.70%  '
Set-TextValue 'D25' '146.29'
Set-TextValue 'E25' '  +0.19%  '
Set-TextValue 'E26' '  -0.92%  '
Set-TextValue 'D27' '7.243'
Set-TextValue 'E27' '  +0.48%  '
Set-TextValue 'D28' '16.11'
Set-TextValue 'E28' '  +0.07%  '
Set-TextValue 'D29' '1.486'
Set-TextValue 'E29' '  -0.92%  '
Set-TextValue 'D30' '0.05671'
Set-TextValue 'E30' '  -3.25%  '
Set-TextValue 'D31' '1.272'
Set-TextValue 'E31' '  -0.85%  '
Set-TextValue 'D32' '3.487'
Set-TextValue 'E32' '  -2.85%  '
Set-TextValue 'D33' '3.362'
Set-TextValue 'E33' '  +2.30%  '
Set-TextValue 'D34' '1.609'
Set-TextValue 'E34' '  +0.10%  '
Set-TextValue 'D35' '2.801'
Set-TextValue 'E35' '  -1.05%  '
Set-TextValue 'D36' '0.9439'
Set-TextValue 'E36' '  -2.61%  '
Set-TextValue 'D37' '2.401'
Set-TextValue 'E37' '  -0.80%  '
Set-TextValue 'D38' '0.5777'
Set-TextValue 'E38' '  -0.61%  '
Set-TextValue 'B39' 'VeChain'
Set-TextValue 'C39' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D39' '0.01600'
Set-TextValue 'E39' '  -0.52%  '
Set-TextValue 'B40' 'FraxShare'
Set-TextValue 'C40' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D40' '5.993'
Set-TextValue 'E40' '  +2.67%  '
Set-TextValue 'D41' '1.071.59'
Set-TextValue 'E41' '  +0.38%  '
Set-TextValue 'D42' '2.581'
Set-TextValue 'E42' '  +0.03%  '
Set-TextValue 'D43' '0.8493'
Set-TextValue 'E43' '  -1.36%  '
Set-TextValue 'E44' '  -0.73%  '
Set-TextValue 'D45' '102.96'
Set-TextValue 'E45' '  -1.22%  '
Set-TextValue 'D46' '1.801.44'
Set-TextValue 'E46' '  -0.18%  '
Set-TextValue 'D47' '58.11'
Set-TextValue 'E47' '  +0.07%  '
Set-TextValue 'E48' '  +2.33%  '
Set-TextValue 'D49' '1.004'
Set-TextValue 'E49' '  -1.31%  '
Set-TextValue 'D50' '0.05314'
Set-TextValue 'E50' '  +2.87%  '
Set-TextValue 'B51' 'EnergySwap'
Set-TextValue 'C51' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D51' '8.049'
Set-TextValue 'E51' '  +0.52%  '
